# Applies the crypto price/volume update described in the commit diff.
# Values in column D (Price) that look like plain numbers must be forced to
# remain text (matching the original inlineStr/shared-string cell type), so we
# temporarily apply a Text number format, assign the value, then restore the
# cell style to Normal so no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '92.110.52'
$ws.Range('E2').Value = '  -2.76%  '
$ws.Range('D3').Value = '3.278.40'
$ws.Range('E3').Value = '  -4.73%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.00'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -5.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '603.54'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.88%  '
$ws.Range('E7').Value = '  -8.01%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.374'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -6.92%  '
$ws.Range('E9').Value = '  +0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.915'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.80%  '
$ws.Range('D11').Value = '3.275.87'
$ws.Range('E11').Value = '  -4.63%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '41.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.84%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.76%  '
$ws.Range('D15').Value = '92.071.84'
$ws.Range('E15').Value = '  -2.71%  '
$ws.Range('D16').Value = '3.889.39'
$ws.Range('E16').Value = '  -5.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000238'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.71%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.90'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -6.47%  '
$ws.Range('D19').Value = '3.278.39'
$ws.Range('E19').Value = '  -4.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '16.89'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.53'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.37'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '480.53'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.435'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -13.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000174'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -9.00%  '
$ws.Range('E26').Value = '  -9.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '88.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.50'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.37%  '
$ws.Range('D29').Value = '3.467.05'
$ws.Range('E29').Value = '  -4.60%  '
$ws.Range('E30').Value = '  +0.11%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '10.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.44%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.134'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -6.63%  '
$ws.Range('E34').Value = '  +0.25%  '
$ws.Range('E35').Value = '  -7.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '27.65'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -10.31%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.516'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -8.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '532.06'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.08%  '
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.16'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.61%  '
$ws.Range('E41').Value = '  -3.45%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.33'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -7.30%  '
$ws.Range('B43').Value = 'ARBITRUM'
$ws.Range('C43').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.841'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -8.11%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '23.81'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('B45').Value = 'ImmutableX'
$ws.Range('C45').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.64'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -2.96%  '
$ws.Range('B46').Value = 'MantraDAO'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0400'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.45%  '
$ws.Range('E48').Value = '  -7.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '51.31'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.72%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.73%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.76'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.28%  '
